$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.82738351924175, 0.8153838604151, 0.88056996714544, 0.81892111811686),
    @(0.86770998809498, 0.81285776754473, 0.87371640094123, 0.83344827633768),
    @(0.8845373489403, 0.82696200048912, 0.88227355876314, 0.87254397919085),
    @(0.88899031765418, 0.81924568196412, 0.89617491395655, 0.8704875235578),
    @(0.90897165768578, 0.83034435195377, 0.88149220963809, 0.84750232386965),
    @(0.89324036302584, 0.84024424199975, 0.87828913846631, 0.87154701734997),
    @(0.88000390920867, 0.81118570577442, 0.88199862548456, 0.87900332263319),
    @(0.89382350518692, 0.84053279048134, 0.88990550290693, 0.87739007628422),
    @(0.88387486895669, 0.86260845777312, 0.8898713420691, 0.88376230656661),
    @(0.86780374744058, 0.82790669460194, 0.87609751653224, 0.86563577545383),
    @(0.87366018903715, 0.82473360730645, 0.88008032699229, 0.8690502179098),
    @(0.88598252369098, 0.85441828256267, 0.89124005269834, 0.8415965991997),
    @(0.87019362970145, 0.79936937037447, 0.87604356286468, 0.85647534002307),
    @(0.8852567096839, 0.83050164791447, 0.87868279899939, 0.84676228046345),
    @(0.87557482896458, 0.8144029109035, 0.8779112034923, 0.86232957467933),
    @(0.88400720488157, 0.81139387453276, 0.87377949960447, 0.84064962183312),
    @(0.86902869913562, 0.80519938234008, 0.87831569256412, 0.84081746470917),
    @(0.87107091843847, 0.83802110030358, 0.87164616302665, 0.85393653214292),
    @(0.87766205322, 0.82018698986832, 0.86773765006325, 0.83815661593337),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    for ($j = 0; $j -lt 4; $j++) {
        $c = $j + 2
        $ws.Cells.Item($r, $c).Value = $data[$i][$j]
    }
}
